# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N on the "Repayment schedule"
# sheet, gives it the same width as column M, moves the selection there,
# and makes that sheet the active tab (matching the tabSelected/activeTab
# change moving from "NewLoanInput" to "Repayment schedule").

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column in front of column N (shifts N->O, O->P, P->Q).
$wsSchedule.Columns("N").Insert() | Out-Null

# Match the width Excel recorded for the newly inserted column (same as
# column M's width).
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Restore the selection on the sheet to where the user left it.
$wsSchedule.Range("S9").Select() | Out-Null

# Make "Repayment schedule" the active sheet/tab (this also clears
# tabSelected on whichever sheet was previously active, i.e. NewLoanInput).
$wsSchedule.Activate()
